$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ColumnHeaders")

# Insert two new rows above row 4 (pushes existing rows 4-12 down to 6-14)
$ws.Rows("4:5").Insert()

# New row 4: ship / New row 5: cruise_ID
# (write column A for both rows first, then column B, so new shared-string
#  entries are appended in the same order Excel recorded them)
$ws.Range("A4").Value = "ship"
$ws.Range("A5").Value = "cruise_ID"
$ws.Range("B4").Value = "Name of vessel from which samples were collected"
$ws.Range("B5").Value = "Cruise identifer for the R/V Tioga"
$ws.Range("C4").Value = "character"
$ws.Range("C5").Value = "character"

# Apply the same wrap-text style used by the rest of column B to the new cells
$ws.Range("B4:B5").WrapText = $true

# Update selection to match the authored state
$ws.Range("C5").Select()
